$wb = $excel.ActiveWorkbook

# --- tc028 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "tc028"
$ws1.Range("A1").Value = "Pagination"
$c = $ws1.Range("A2")
$c.NumberFormat = "@"
$c.Value = "1"
$c.ClearFormats()

# --- tc029 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "tc029"
$ws2.Range("A1").Value = "Pagination"
$d = $ws2.Range("A2")
$d.NumberFormat = "@"
$d.Value = "2"
$d.ClearFormats()

$wb.Worksheets.Item(1).Activate()
